$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '25.849.45'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.635.70'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +0.74%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '214.97'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.10%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5084'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06426'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.90%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '20.40'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +5.53%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07798'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.32%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '4.267'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +1.41%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.861.36'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.79%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '1.634.93'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.69%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.5604'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.78%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0₅7669'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +2.55%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '63.21'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.23%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '25.857.65'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.35%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '1.004'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.17%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '193.56'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.13%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.390'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.18%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '9.952'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.98%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '6.156'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +2.90%  '
$ws.Range('E24').Value = '  +0.14%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.793'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -4.24%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '138.63'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -1.92%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.1233'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.41%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.853'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +2.36%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.04965'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +2.29%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.306'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +2.46%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.251'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +3.20%  '
$ws.Range('E34').Value = '  +2.41%  '
$ws.Range('E35').Value = '  +0.77%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.9033'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +1.32%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.578'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +1.66%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.5563'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.33%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.134.69'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +2.06%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.01568'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +1.70%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.9966'
$ws.Range('D41').Style = "Normal"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '99.22'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +2.27%  '
$ws.Range('E43').Value = '  -0.89%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.7998'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.57%  '
$ws.Range('E45').Value = '  -1.49%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '55.52'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.92%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.4262'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -3.65%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '7.757'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +2.91%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.05054'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -1.32%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.9993'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.41%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.002'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.05%  '
